# Update "disability_prevalence.xlsx" (Kobuleti) to the new, revised layout:
#  - New title text referring to the "Unified database of targeted social
#    assistance program" (row 1, merged across A1:I1).
#  - A new data row "family with disabilities Persons" inserted above the
#    existing data row, which itself is relabeled "disabilities Persons"
#    and gets a fresh set of figures.
#  - The source note moves down one row (now row 6) and keeps its border
#    and merge, spanning the (now) 6-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row above the old data row (old row 4) so we end up
#    with two stacked data rows instead of one. Excel shifts the old
#    row 4 ("Number of disability persons" + values) down to row 5, and
#    the old row 5 (Source, merged A5:H5) down to row 6, along with all
#    formatting/merges.
# ---------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------
# 2. Title (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Kobuleti Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:I1").VerticalAlignment = -4108     # xlCenter
$ws.Range("A1:I1").WrapText = $true
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Bold = $true
$ws.Rows.Item(1).RowHeight = 51
$ws.Range("A1:I1").Interior.Pattern = -4142      # xlNone

# Row 2 reverts to the sheet's (new) default height - no custom height
$ws.Rows.Item(2).RowHeight = 14.5

# Column A narrows a bit under the new layout; other columns drop their
# explicit custom width and fall back to the sheet default.
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns("B:P").ColumnWidth = 8.43

# ---------------------------------------------------------------------
# 3. Row 4 - new data row "family with disabilities Persons "
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.Bold = $false
$ws.Range("A4").HorizontalAlignment = -4131      # xlLeft
$ws.Range("A4").VerticalAlignment = -4108        # xlCenter
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Interior.Pattern = 1             # xlSolid
$ws.Range("A4").Interior.ThemeColor = 1          # white (background 1)
$ws.Range("A4").Borders.Item(9).LineStyle = 1    # xlEdgeTop continuous
$ws.Range("A4").Borders.Item(9).Weight = 2       # xlThin
$ws.Rows.Item(4).RowHeight = 24.75

$rowB4 = @(1653,1622,1573,1700,1789,1845,1863,1936)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $rowB4[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Interior.Pattern = 1
    $cell.Interior.ThemeColor = 1
}

# ---------------------------------------------------------------------
# 4. Row 5 - relabeled data row "disabilities Persons " with new figures
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.Bold = $false
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Interior.Pattern = 1
$ws.Range("A5").Interior.ThemeColor = 1
$ws.Range("A5").Borders.Item(9).LineStyle = 0    # clear old top border
$ws.Range("A5").Borders.Item(12).LineStyle = 1   # xlEdgeBottom continuous
$ws.Range("A5").Borders.Item(12).Weight = 2
$ws.Rows.Item(5).RowHeight = 21

$rowB5 = @(1932,1907,1844,1986,2067,2131,2136,2206)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $rowB5[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Interior.Pattern = 1
    $cell.Interior.ThemeColor = 1
    if ($cols[$i] -eq "I") {
        $cell.Borders.Item(12).LineStyle = 1
        $cell.Borders.Item(12).Weight = 2
    }
}

# ---------------------------------------------------------------------
# 5. Row 6 - Source note, now spans A6:H6 (shifted automatically by the
#    row insert, but make sure the merge / sizing matches explicitly).
# ---------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 27.75
$ws.Range("A6:H6").Merge()

$wb.Save()
